$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3266.5334
$ws.Range("J2").Value = 3649.6667
$ws.Range("L2").Value = 3649.6667
$ws.Range("N2").Value = -3875.6667
$ws.Range("H6").Value = 198.27272
$ws.Range("I6").Value = 166.375
$ws.Range("K6").Value = 499.125
$ws.Range("M6").Value = -387.125
$ws.Range("H7").Value = 12084.333
$ws.Range("I7").Value = 8249.5
$ws.Range("J7").Value = 14001.75
$ws.Range("K7").Value = 8249.5
$ws.Range("L7").Value = 14001.75
$ws.Range("M7").Value = -8137.5
$ws.Range("N7").Value = -14225.75
$ws.Range("H9").Value = 167
$ws.Range("J9").Value = 167
$ws.Range("L9").Value = 167
$ws.Range("N9").Value = -505
$ws.Range("H10").Value = 12667.667
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 12667.667
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 12667.667
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -13253.667
$ws.Range("H14").Value = 12084.333
$ws.Range("I14").Value = 8249.5
$ws.Range("J14").Value = 14001.75
$ws.Range("K14").Value = 8249.5
$ws.Range("L14").Value = 14001.75
$ws.Range("M14").Value = -8058.5
$ws.Range("N14").Value = -14383.75
$ws.Range("H16").Value = 6141.6
$ws.Range("I16").Value = 350
$ws.Range("J16").Value = 10002.667
$ws.Range("K16").Value = 350
$ws.Range("L16").Value = 10002.667
$ws.Range("M16").Value = -120
$ws.Range("N16").Value = -10462.667
$ws.Range("H40").Value = 1500
$ws.Range("J40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("N40").Value = -1850
$ws.Range("H98").Value = 1327.2354
$ws.Range("I98").Value = 1327.2354
$ws.Range("K98").Value = 1327.2354
$ws.Range("M98").Value = 170.7646
$ws.Range("H122").Value = 1327.2354
$ws.Range("I122").Value = 1327.2354
$ws.Range("K122").Value = 3981.7062
$ws.Range("M122").Value = -1531.7062
$ws.Range("H135").Value = 235
$ws.Range("I135").Value = 235
$ws.Range("K135").Value = 2115
$ws.Range("M135").Value = 420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 464.57144
$ws.Range("I5").Value = 500.4
$ws.Range("J5").Value = 375
$ws.Range("K5").Value = 500.4
$ws.Range("L5").Value = 375
$ws.Range("M5").Value = -388.4
$ws.Range("N5").Value = -599
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 5051
$ws.Range("I12").Value = 350
$ws.Range("J12").Value = 9752
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 9752
$ws.Range("M12").Value = -177
$ws.Range("N12").Value = -10098
$ws.Range("H13").Value = 17373.5
$ws.Range("J13").Value = 18164.666
$ws.Range("L13").Value = 18164.666
$ws.Range("N13").Value = -18452.666
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 1900
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 8672.666999999999
$ws.Range("J22").Value = 12509
$ws.Range("L22").Value = 12509
$ws.Range("N22").Value = -13107
$ws.Range("H25").Value = 4605.4
$ws.Range("I25").Value = 752.25
$ws.Range("K25").Value = 752.25
$ws.Range("M25").Value = -350.25
$ws.Range("H30").Value = 8190
$ws.Range("I30").Value = 916.6667
$ws.Range("J30").Value = 30010
$ws.Range("K30").Value = 916.6667
$ws.Range("L30").Value = 30010
$ws.Range("M30").Value = -766.6667
$ws.Range("N30").Value = -30310
$ws.Range("H36").Value = 2747.75
$ws.Range("I36").Value = 2747.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2747.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2401.75
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 8834.833000000001
$ws.Range("I61").Value = 8601.799999999999
$ws.Range("K61").Value = 8601.799999999999
$ws.Range("M61").Value = -8389.799999999999
$ws.Range("H92").Value = 110550
$ws.Range("J92").Value = 110550
$ws.Range("L92").Value = 110550
$ws.Range("N92").Value = -115542
$ws.Range("H122").Value = 2998.5
$ws.Range("I122").Value = 2998.5
$ws.Range("K122").Value = 8995.5
$ws.Range("M122").Value = -6545.5
$ws.Range("H136").Value = 8834.833000000001
$ws.Range("I136").Value = 8601.799999999999
$ws.Range("K136").Value = 25805.4
$ws.Range("M136").Value = -23255.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 464.57144
$ws.Range("I4").Value = 500.4
$ws.Range("J4").Value = 375
$ws.Range("K4").Value = 500.4
$ws.Range("L4").Value = 375
$ws.Range("M4").Value = -385.4
$ws.Range("N4").Value = -605
$ws.Range("H10").Value = 679
$ws.Range("I10").Value = 358
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 358
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -218
$ws.Range("N10").Value = -1280
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253
$ws.Range("H134").Value = 9423.611000000001
$ws.Range("I134").Value = 8213.235000000001
$ws.Range("K134").Value = 24639.705
$ws.Range("M134").Value = -22104.705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 436.66666
$ws.Range("I11").Value = 985
$ws.Range("J11").Value = 162.5
$ws.Range("K11").Value = 985
$ws.Range("L11").Value = 162.5
$ws.Range("M11").Value = -845
$ws.Range("N11").Value = -442.5
$ws.Range("H62").Value = 5949.5
$ws.Range("I62").Value = 5899
$ws.Range("K62").Value = 5899
$ws.Range("M62").Value = -5275
$ws.Range("H65").Value = 5949.5
$ws.Range("I65").Value = 5899
$ws.Range("K65").Value = 29495
$ws.Range("M65").Value = -26375
$ws.Range("H95").Value = 24357.4
$ws.Range("J95").Value = 24357.4
$ws.Range("L95").Value = 24357.4
$ws.Range("N95").Value = -29849.4
$ws.Range("H99").Value = 4962.6
$ws.Range("I99").Value = 6933
$ws.Range("J99").Value = 2007
$ws.Range("K99").Value = 6933
$ws.Range("L99").Value = 2007
$ws.Range("M99").Value = -5435
$ws.Range("N99").Value = -5003
$ws.Range("H122").Value = 1689.8
$ws.Range("I122").Value = 1653.8334
$ws.Range("J122").Value = 1743.75
$ws.Range("K122").Value = 4961.5002
$ws.Range("L122").Value = 5231.25
$ws.Range("M122").Value = -2511.5002
$ws.Range("N122").Value = -10131.25
$ws.Range("H126").Value = 4962.6
$ws.Range("I126").Value = 6933
$ws.Range("J126").Value = 2007
$ws.Range("K126").Value = 20799
$ws.Range("L126").Value = 6021
$ws.Range("M126").Value = -18329
$ws.Range("N126").Value = -10961
$ws.Range("H132").Value = 1771.1538
$ws.Range("I132").Value = 1364.6364
$ws.Range("K132").Value = 4093.9092
$ws.Range("M132").Value = -1563.9092
$ws.Range("H134").Value = 3857.6155
$ws.Range("I134").Value = 3761.25
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 11283.75
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -8748.75
$ws.Range("N134").Value = -20112
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H139").Value = 2080.2727
$ws.Range("I139").Value = 2080.2727
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 6240.8181
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -1100.8181
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29000
$ws.Range("J15").Value = 29000
$ws.Range("L15").Value = 29000
$ws.Range("N15").Value = -29576
$ws.Range("H81").Value = 29000
$ws.Range("J81").Value = 29000
$ws.Range("L81").Value = 29000
$ws.Range("N81").Value = -30996
$ws.Range("H84").Value = 29000
$ws.Range("J84").Value = 29000
$ws.Range("L84").Value = 87000
$ws.Range("N84").Value = -96984
$ws.Range("H97").Value = 630.9091
$ws.Range("I97").Value = 541.125
$ws.Range("J97").Value = 870.3333
$ws.Range("K97").Value = 541.125
$ws.Range("L97").Value = 870.3333
$ws.Range("M97").Value = -45.125
$ws.Range("N97").Value = -1862.3333
$ws.Range("H126").Value = 1912
$ws.Range("I126").Value = 1912
$ws.Range("K126").Value = 5736
$ws.Range("M126").Value = -3266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8062.5
$ws.Range("I132").Value = 6125
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 18375
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -15845
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 1500
$ws.Range("M132").Value = 1030
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
